# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-22) needs to be re-sorted in ascending
# order by the "Periodo Mora" column (E). The "Valor Mora" column (F) must
# travel together with its corresponding period so that each period keeps
# its original value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 22

# Read the current Periodo Mora (E) and Valor Mora (F) values for every row
# in the table.
$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodo = $ws.Cells.Item($r, 5).Value2
    $valor = $ws.Cells.Item($r, 6).Value2
    $rows += , @($periodo, $valor)
}

# Sort the captured rows by Periodo Mora ascending (the values are 4-digit
# text/numeric codes such as 2308, 2309, ... 2402, so a plain string/number
# sort gives the desired chronological order).
$sorted = $rows | Sort-Object { [string]$_[0] }

# Write the sorted Periodo Mora / Valor Mora pairs back into the range,
# leaving every other column (Tipo Doc, N Doc, Nombre, Salario Basico,
# Novedad de Ingreso/Retiro, Observaciones) untouched since they are
# identical for every row of this worker.
for ($idx = 0; $idx -lt $sorted.Count; $idx++) {
    $r = $firstRow + $idx
    $pair = $sorted[$idx]
    $ws.Cells.Item($r, 5).Value2 = $pair[0]
    $ws.Cells.Item($r, 6).Value2 = $pair[1]
}
